$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new value would otherwise be
# auto-detected as a number by Excel (single-decimal numeric strings),
# so they are stored as text exactly like the rest of the Price column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "67.645.37"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.641.37"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "605.14"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "154.58"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").Value = "2.636.30"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("E10").Value = "  +7.46%  "
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").Value = "28.05"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "3.117.18"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "67.595.61"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "2.638.38"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "11.33"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "366.22"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "7.63"
$ws.Range("E21").Value = "  -4.47%  "
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("E23").Value = "  +7.45%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").Value = "66.16"
$ws.Range("E26").Value = "  -8.03%  "
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "583.37"
$ws.Range("E29").Value = "  -7.27%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  -3.36%  "
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").Value = "158.06"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("D39").Value = "19.51"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").Value = "5.34"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "156.50"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("D50").Value = "21.02"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("E51").Value = "  +0.44%  "
